# Update countries & provincias Spain
# - Swap ranking order of Italia/Irak (Italia now ahead of Irak)
# - Swap ranking order of Moldavia/Armenia (Moldavia now ahead of Armenia)
# - Refresh the "Datos actualizados" timestamp
# - Update the underlying case/death statistics that changed with the new snapshot

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner (row 1)
$ws.Cells.Item(1,1).Value = "Datos actualizados a 21 de Octubre de 2020 a las 17:44"

# --- Row 18 / Row 19: Italia moves ahead of Irak ---
$ws.Cells.Item(18,1).Value = "Italia"
$ws.Cells.Item(18,2).Value = 449648
$ws.Cells.Item(18,3).Value = 15199
$ws.Cells.Item(18,4).Value = 257374
$ws.Cells.Item(18,5).Value = 155442
$ws.Cells.Item(18,7).Value = 127
$ws.Cells.Item(18,8).Value = 36832

$ws.Cells.Item(19,1).Value = "Irak"
$ws.Cells.Item(19,2).Value = 438265
$ws.Cells.Item(19,3).Value = 3667
$ws.Cells.Item(19,4).Value = 369010
$ws.Cells.Item(19,5).Value = 58837
$ws.Cells.Item(19,7).Value = 52
$ws.Cells.Item(19,8).Value = 10418

# --- Row 60 / Row 61: Moldavia moves ahead of Armenia ---
$ws.Cells.Item(60,1).Value = "Moldavia"
$ws.Cells.Item(60,2).Value = 68791
$ws.Cells.Item(60,3).Value = 833
$ws.Cells.Item(60,4).Value = 49702
$ws.Cells.Item(60,5).Value = 17459
$ws.Cells.Item(60,7).Value = 13
$ws.Cells.Item(60,8).Value = 1630

$ws.Cells.Item(61,1).Value = "Armenia"
$ws.Cells.Item(61,2).Value = 68530
$ws.Cells.Item(61,3).Value = 1836
$ws.Cells.Item(61,4).Value = 49219
$ws.Cells.Item(61,5).Value = 18190
$ws.Cells.Item(61,7).Value = 20
$ws.Cells.Item(61,8).Value = 1121

# --- Remaining statistic refreshes (country identities unchanged) ---

# Row 4: Estados Unidos
$ws.Cells.Item(4,2).Value = 8524940
$ws.Cells.Item(4,3).Value = 3990
$ws.Cells.Item(4,5).Value = 2750592
$ws.Cells.Item(4,7).Value = 110
$ws.Cells.Item(4,8).Value = 226294

# Row 5: India
$ws.Cells.Item(5,2).Value = 7670537
$ws.Cells.Item(5,3).Value = 21379
$ws.Cells.Item(5,4).Value = 6812889
$ws.Cells.Item(5,5).Value = 741580
$ws.Cells.Item(5,7).Value = 118
$ws.Cells.Item(5,8).Value = 116068

# Row 17: Chile
$ws.Cells.Item(17,2).Value = 495637
$ws.Cells.Item(17,3).Value = 1159
$ws.Cells.Item(17,4).Value = 468269
$ws.Cells.Item(17,5).Value = 13649
$ws.Cells.Item(17,7).Value = 17
$ws.Cells.Item(17,8).Value = 13719

# Row 31: Canada
$ws.Cells.Item(31,2).Value = 204479
$ws.Cells.Item(31,3).Value = 791
$ws.Cells.Item(31,4).Value = 172462
$ws.Cells.Item(31,5).Value = 22214
$ws.Cells.Item(31,7).Value = 9
$ws.Cells.Item(31,8).Value = 9803

# Row 41: Republica Dominicana
$ws.Cells.Item(41,2).Value = 122398
$ws.Cells.Item(41,3).Value = 425
$ws.Cells.Item(41,4).Value = 100051
$ws.Cells.Item(41,5).Value = 20141
$ws.Cells.Item(41,7).Value = 2
$ws.Cells.Item(41,8).Value = 2206

# Row 51: Japon
$ws.Cells.Item(51,2).Value = 93933
$ws.Cells.Item(51,3).Value = 453
$ws.Cells.Item(51,4).Value = 87107
$ws.Cells.Item(51,5).Value = 5147
$ws.Cells.Item(51,7).Value = 3
$ws.Cells.Item(51,8).Value = 1679

# Row 52: Suiza
$ws.Cells.Item(52,5).Value = 35134
$ws.Cells.Item(52,7).Value = 7
$ws.Cells.Item(52,8).Value = 2029

# Row 95: Albania
$ws.Cells.Item(95,2).Value = 17948
$ws.Cells.Item(95,3).Value = 297
$ws.Cells.Item(95,4).Value = 10341
$ws.Cells.Item(95,5).Value = 7145
$ws.Cells.Item(95,7).Value = 4
$ws.Cells.Item(95,8).Value = 462

# Row 115: Jamaica
$ws.Cells.Item(115,2).Value = 8445
$ws.Cells.Item(115,3).Value = 71
$ws.Cells.Item(115,4).Value = 4016
$ws.Cells.Item(115,5).Value = 4255

# Row 122: Cuba
$ws.Cells.Item(122,2).Value = 6368
$ws.Cells.Item(122,3).Value = 63
$ws.Cells.Item(122,4).Value = 5814
$ws.Cells.Item(122,5).Value = 427

# Row 123: Bahamas
$ws.Cells.Item(123,2).Value = 6051
$ws.Cells.Item(123,3).Value = 128
$ws.Cells.Item(123,4).Value = 3633
$ws.Cells.Item(123,5).Value = 2291
$ws.Cells.Item(123,7).Value = 3
$ws.Cells.Item(123,8).Value = 127

# Row 170: San Marino
$ws.Cells.Item(170,2).Value = 774
$ws.Cells.Item(170,3).Value = 8
$ws.Cells.Item(170,4).Value = 690
$ws.Cells.Item(170,5).Value = 42
